# Populate the newly-added "-" (not available) placeholder cells in
# columns I (Unemployment rate), J (Poverty rate) and K (HDI) for the
# country/year rows where that data was missing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where I, J and K were all empty and now all get "-"
$fullRows = @(2,3,4,5,6,7,8,9,14,15,16,17)
foreach ($r in $fullRows) {
    foreach ($col in @("I","J","K")) {
        $ws.Range("$col$r").Value = "-"
    }
}

# Rows where only J (poverty rate) was empty; I and K already had values
$partialRows = @(25,30,31,32)
foreach ($r in $partialRows) {
    $ws.Range("J$r").Value = "-"
}

# Restore the selection to the single cell E5 (scrolled back to A1)
$null = $ws.Range("E5").Select()
